$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 35, pushing the existing rows 35-41 down to 36-42
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly data point
$ws.Range("A35").Value = 8
$ws.Range("B35").Value = "Terminal La Palmera de La Serena"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 45180
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112013
$ws.Range("G35").Value = "Alcachofa"
$ws.Range("H35").Value = "Española"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 9500
$ws.Range("N35").Value = "$/caja 30 unidades"
$ws.Range("O35").Value = "Provincia del Elquí"
$ws.Range("P35").Value = 317
$ws.Range("Q35").Value = 30
$ws.Range("R35").Value = "Hortaliza"
